$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 5 (A5 already has "C - Otoshidama ABC085") ---
# D5 re-uses the SAME note text already stored for row 4 ("X + Y + Z = ...")
# -- copy the existing cell's value across so it lands on the same shared
#    string slot instead of minting a duplicate.
$ws.Range("B5").Value = "解けた"
$ws.Range("D5").Value = $ws.Range("D4").Value2

# --- New rows 6-12: problem name / solved flag / note ---
$ws.Range("A6").Value  = "C - たくさんの数式 / Many Formulas ABC061"
$ws.Range("B6").Value  = "解けた"
$ws.Range("D6").Value  = "bit全探索"

$ws.Range("A7").Value  = "C - Train Ticket ABC79"
$ws.Range("B7").Value  = "解けた"
$ws.Range("D7").Value  = "bit全探索"

$ws.Range("A8").Value  = "C - All Green ABC104"
$ws.Range("B8").Value  = "解けた"
$ws.Range("D8").Value  = "dp[i][j] : i番目までの問題をj問解いた時の最大獲得点数　で動的計画法"

$ws.Range("A9").Value  = "A - 高橋君とお肉 ARC029"
$ws.Range("B9").Value  = "解けた"
$ws.Range("D9").Value  = "bit全探索で１が立っていれば肉焼き機Aで０であればBで焼いてその最小時間を出力"

$ws.Range("A10").Value = "D - 派閥"
$ws.Range("B10").Value = "解けた"
$ws.Range("D10").Value = "bit全探索で１が立っている者同士を同じ派閥にして条件を満たしているかを確認"

$ws.Range("A11").Value = "A - 深さ優先探索 ATC001"
$ws.Range("B11").Value = "解けた"
$ws.Range("D11").Value = "再帰関数を用いたdfsとスタックを用いたdfsを学習"

$ws.Range("A12").Value = "B - 埋め立て ARC031"
$ws.Range("B12").Value = "解けた"
$ws.Range("D12").Value = "スタックを用いたdfs　埋め立て地の候補は高々１００個だから全探索"

# --- Dates (column C), rows 5-12: copy the already-date-formatted C4 cell's
#     format down so every row shares the SAME style slot instead of each
#     NumberFormat assignment minting a brand-new (duplicate) style ---
$ws.Range("C4").Copy()
$ws.Range("C5:C12").PasteSpecial(-4122)
$ws.Range("C5:C12").Value = 43409
$excel.CutCopyMode = 0

# --- Column widths (characters) ---
$ws.Columns.Item(1).ColumnWidth = 44.857142857142854
$ws.Columns.Item(2).ColumnWidth = 16.714285714285715
$ws.Columns.Item(3).ColumnWidth = 13.428571428571429
$ws.Columns.Item(4).ColumnWidth = 68.42857142857143
$ws.Columns.Item(5).ColumnWidth = 61.0

# --- Selection moved from A5 to D12 ---
$ws.Range("D12").Select()

Write-Output "edit applied"
